# The underlying smoothing data (column B, rows 1-100) on the
# "output2SaltedSmoothened" sheet was regenerated with a corrected
# smoothing routine. Column A (the X values, 0..99) is untouched;
# only the smoothed Y values in column B change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    81.6666666666666,
    158.25,
    64,
    203.2,
    114.4,
    52.2,
    -5.8,
    37,
    113.8,
    310.39999999999998,
    475.6,
    689.2,
    836.2,
    921.8,
    942.2,
    1082.8,
    1101.8,
    1344.8,
    1399.8,
    1571.6,
    1741,
    2021.4,
    2086.8000000000002,
    2334,
    2468.8000000000002,
    2676.6,
    2822.2,
    3154.2,
    3442.6,
    3806.4,
    4003.2,
    4137.3999999999996,
    4238,
    4364,
    4640.3999999999996,
    4955.2,
    5367.4,
    5805.6,
    6083.4,
    6332.8,
    6626,
    6901,
    7104,
    7626.4,
    7939.6,
    8171.2,
    8613.4,
    9053,
    9398.6,
    9926.6,
    10384.799999999999,
    10802,
    11160.6,
    11452.6,
    11856.4,
    12252.2,
    12550.6,
    13171.8,
    13682.4,
    14158.2,
    14673.8,
    15282.4,
    15789.8,
    16232,
    16609.2,
    17245.599999999999,
    17682.8,
    18220.2,
    18945.8,
    19677,
    20112.400000000001,
    20789.599999999999,
    21217.599999999999,
    21759,
    22225,
    22786.2,
    23427.4,
    24024.400000000001,
    24693.8,
    25279.4,
    25958.799999999999,
    26557.8,
    27379,
    28002.799999999999,
    28833,
    29612.799999999999,
    30332.6,
    31026.6,
    31661.200000000001,
    32221.8,
    32875.800000000003,
    33625.599999999999,
    34337,
    35173,
    36030.199999999997,
    36743,
    37380.6,
    38056.199999999997,
    38340.25,
    38631.666666666599
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
